$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: copy formatting from E1 (bold header style) then set text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F6 with the recorded timestamps
$ws.Range("F2").Value = "2021-10-05 13:42:26.832164"
$ws.Range("F3").Value = "2021-10-05 13:42:26.832174"
$ws.Range("F4").Value = "2021-10-05 13:42:26.832177"
$ws.Range("F5").Value = "2021-10-05 13:42:26.832180"
$ws.Range("F6").Value = "2021-10-05 13:42:26.832183"

$excel.CutCopyMode = $false
